$wb = $excel.ActiveWorkbook

# Sheet "展览" updates (column F = 想去人数)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 3643
$ws1.Range("F5").Value = 2229
$ws1.Range("F8").Value = 176
$ws1.Range("F10").Value = 74
$ws1.Range("F11").Value = 1341
$ws1.Range("F13").Value = 2001

# Sheet "全部类型" updates (column F = 想去人数)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 3643
$ws4.Range("F5").Value = 2229
$ws4.Range("F9").Value = 176
$ws4.Range("F11").Value = 74
$ws4.Range("F14").Value = 1341
$ws4.Range("F16").Value = 2001
